$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- Paragraph 1: "Liaison radio" -> "Liaison radio – chercher bon module" (+ fr-FR lang) ---
$p1 = $d.Paragraphs.Item(1)
$frag1 = '<w:p ' + $ns + ' w14:paraId="487DA29F" w14:textId="2C36B4A6" w:rsidR="00E718C5" w:rsidRDefault="00D61F8E" w:rsidP="00D61F8E">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
  '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-FR"/></w:rPr><w:t>Liaison radio</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> &#8211; chercher bon m</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-FR"/></w:rPr><w:t>odule</w:t></w:r>' + `
  '</w:p>'
$p1.Range.InsertXML($frag1)

# --- Paragraph 2: "Acceleromètre :" -- drop proofErr, merge into one run ---
$p2 = $d.Paragraphs.Item(2)
$frag2 = '<w:p ' + $ns + ' w14:paraId="49BD1860" w14:textId="1CF33EDB" w:rsidR="00D61F8E" w:rsidRDefault="30C9928A" w:rsidP="00D61F8E">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
  '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Acceleromètre :</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> MPU6050</w:t></w:r>' + `
  '</w:p>'
$p2.Range.InsertXML($frag2)

# --- Paragraph 4: "Module peltier ? Panneau solaire ?" -- drop proofErr, merge runs ---
$p4 = $d.Paragraphs.Item(4)
$frag4 = '<w:p ' + $ns + ' w14:paraId="29F92FFF" w14:textId="64A2A1C7" w:rsidR="00D61F8E" w:rsidRDefault="00D61F8E" w:rsidP="00D61F8E">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
  '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Module peltier ?</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Panneau solaire ?</w:t></w:r>' + `
  '</w:p>'
$p4.Range.InsertXML($frag4)

Write-Host "Edits applied."
